$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.585.24'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '3.026.06'
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.42%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.528'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.63%  '
$ws.Range("D9").Value = '3.025.07'
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("E10").Value = '  -3.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.90'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.453'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000231'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.69%  '
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("D16").Value = '3.524.19'
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").Value = '62.458.99'
$ws.Range("E18").Value = '  -1.22%  '
$ws.Range("D19").Value = '3.022.31'
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '465.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.689'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("E25").Value = '  -5.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -1.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.64'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.07%  '
$ws.Range("E34").Value = '  -1.83%  '
$ws.Range("D35").Value = '0.0₃0811'
$ws.Range("E35").Value = '  -1.04%  '
$ws.Range("E36").Value = '  -2.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.52'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.16'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.114'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '397.06'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.45%  '
$ws.Range("E44").Value = '  -3.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0359'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").Value = '2.762.33'
$ws.Range("E46").Value = '  -2.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '37.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.110'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.16%  '
